$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text so that numeric-looking
# strings (e.g. "1.006", "2.350") keep their exact characters instead of
# being auto-coerced into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.138.34'
$ws.Range("E2").Value = '  +1.93%  '
$ws.Range("D3").Value = '1.718.59'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '319.18'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").Value = '0.3972'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.4125'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("D9").Value = '1.529'
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("D10").Value = '1.005'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").Value = '53.56'
$ws.Range("E11").Value = '  +4.93%  '
$ws.Range("D12").Value = '0.08941'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '7.748'
$ws.Range("E13").Value = '  +7.69%  '
$ws.Range("D14").Value = '25.05'
$ws.Range("E14").Value = '  +6.83%  '
$ws.Range("D15").Value = '8.175'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '0.00001385'
$ws.Range("E16").Value = '  +4.17%  '
$ws.Range("D17").Value = '1.712.85'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '100.92'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").Value = '0.07156'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = '20.28'
$ws.Range("E20").Value = '  +2.94%  '
$ws.Range("D21").Value = '7.495'
$ws.Range("E21").Value = '  +6.36%  '
$ws.Range("D22").Value = '1.007'
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").Value = '14.55'
$ws.Range("E23").Value = '  +2.37%  '
$ws.Range("D24").Value = '25.115.65'
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").Value = '3.120'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '2.351'
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '23.12'
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '165.20'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '8.856'
$ws.Range("E29").Value = '  +19.78%  '
$ws.Range("D30").Value = '139.87'
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("D31").Value = '5.218'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Value = '7.814'
$ws.Range("E32").Value = '  +9.77%  '
$ws.Range("D33").Value = '1.902.00'
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("D34").Value = '0.09016'
$ws.Range("E34").Value = '  +4.82%  '
$ws.Range("D35").Value = '1.081'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02989'
$ws.Range("E36").Value = '  +9.55%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2799'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.974'
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '11.07'
$ws.Range("E39").Value = '  -3.68%  '
$ws.Range("D40").Value = '14.69'
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.09286'
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.8173'
$ws.Range("E42").Value = '  +6.42%  '
$ws.Range("D43").Value = '1.483'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").Value = '16.70'
$ws.Range("E44").Value = '  +4.29%  '
$ws.Range("D45").Value = '0.7418'
$ws.Range("E45").Value = '  +3.32%  '
$ws.Range("D46").Value = '2.646'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").Value = '4.294'
$ws.Range("E47").Value = '  +1.56%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = '1.353'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").Value = '140.73'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '93.21'
$ws.Range("E51").Value = '  +4.29%  '
